$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "56.504.35"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").Value = "2.380.12"
$ws.Range("E3").Value = "  -2.91%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "503.33"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("D6").Value = "130.67"
$ws.Range("E6").Value = "  -1.89%  "
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("E8").Value = "  -2.32%  "
$ws.Range("D9").Value = "2.386.97"
$ws.Range("E9").Value = "  -2.62%  "
$ws.Range("D10").Value = "0.0987"
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "0.326"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").Value = "4.73"
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("D14").Value = "2.803.82"
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("D15").Value = "56.448.51"
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("D16").Value = "21.69"
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").Value = "2.354.35"
$ws.Range("E18").Value = "  -4.69%  "
$ws.Range("E19").Value = "  -2.62%  "
$ws.Range("E20").Value = "  -2.07%  "
$ws.Range("D21").Value = "308.28"
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("D22").Value = "6.26"
$ws.Range("E22").Value = "  -2.22%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "65.54"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").Value = "0.997"
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("E26").Value = "  -3.82%  "
$ws.Range("D27").Value = "0.149"
$ws.Range("E27").Value = "  -3.52%  "
$ws.Range("D28").Value = "7.33"
$ws.Range("E28").Value = "  -3.64%  "
$ws.Range("D29").Value = "172.29"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("D30").Value = "0.0₃0718"
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("E31").Value = "  -2.96%  "
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("E33").Value = "  -6.64%  "
$ws.Range("E34").Value = "  -4.10%  "
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("E36").Value = "  -1.88%  "
$ws.Range("D37").Value = "1.18"
$ws.Range("E37").Value = "  -5.24%  "
$ws.Range("E38").Value = "  -1.75%  "
$ws.Range("D39").Value = "36.13"
$ws.Range("E39").Value = "  -1.16%  "
$ws.Range("D40").Value = "0.799"
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("E41").Value = "  -4.49%  "
$ws.Range("D42").Value = "131.14"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("E44").Value = "  -3.36%  "
$ws.Range("D45").Value = "0.566"
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("D47").Value = "242.01"
$ws.Range("E47").Value = "  -5.79%  "
$ws.Range("D48").Value = "0.0484"
$ws.Range("E48").Value = "  -2.22%  "
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("D50").Value = "17.18"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("E51").Value = "  -2.75%  "

Write-Host "Updated 77 cells (21 with text format applied)"
